$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy formatting from the prior pair of rows (rows 2 and 3, same station order)
# onto the two new rows so the new cells inherit the same styles (s="1","3","2").
$ws.Range("A2:F3").Copy()
$ws.Range("A34:F35").PasteSpecial(-4122)

# Row 34: 四方坪站
$ws.Range("A34").Value = 45978
$ws.Range("B34").Value = "四方坪站"
$ws.Range("C34").Value = 9597.15
$ws.Range("D34").Value = 8627.19
$ws.Range("E34").Value = 3109.43
$ws.Range("F34").Value = 397

# Row 35: 高岭站
$ws.Range("A35").Value = 45978
$ws.Range("B35").Value = "高岭站"
$ws.Range("C35").Value = 4848.63
$ws.Range("D35").Value = 4379.17
$ws.Range("E35").Value = 1240.82
$ws.Range("F35").Value = 175

# Match the updated active selection shown in the target workbook.
$ws.Range("H34").Select()
